$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I24").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("I27").Value = 0

$ws.Range("I28").Select()
